# This edit refreshes the weekly 'Fruta / hortaliza' price report for
# 'Vega Modelo de Temuco - Tuna': two brand-new observations are inserted
# at rows 7-8 (pushing the previously-existing rows 7-70 down by two), and
# the sheet grows by two rows overall (old rows 69-70 become rows 71-72).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the full replacement block for columns D:T, rows 7-72, in one shot.
$data = New-Object 'object[,]' 66,17
$data[0,0] = 44970; $data[0,1] = 9; $data[0,2] = 'Fruta'; $data[0,3] = 100107; $data[0,4] = 'Otros'; $data[0,5] = 100107011; $data[0,6] = 'Tuna'; $data[0,7] = 'Sin especificar'; $data[0,8] = 'Primera'; $data[0,9] = 35; $data[0,10] = 20000; $data[0,11] = 20000; $data[0,12] = 20000; $data[0,13] = '$/caja 18 kilos'; $data[0,14] = 'Provincia de Los Andes'; $data[0,15] = 1111; $data[0,16] = 18
$data[1,0] = 44970; $data[1,1] = 9; $data[1,2] = 'Fruta'; $data[1,3] = 100107; $data[1,4] = 'Otros'; $data[1,5] = 100107011; $data[1,6] = 'Tuna'; $data[1,7] = 'Sin especificar'; $data[1,8] = 'Segunda'; $data[1,9] = 40; $data[1,10] = 15000; $data[1,11] = 15000; $data[1,12] = 15000; $data[1,13] = '$/caja 18 kilos'; $data[1,14] = 'Provincia de Los Andes'; $data[1,15] = 833; $data[1,16] = 18
$data[2,0] = 44321; $data[2,1] = 9; $data[2,2] = 'Fruta'; $data[2,3] = 100107; $data[2,4] = 'Otros'; $data[2,5] = 100107011; $data[2,6] = 'Tuna'; $data[2,7] = 'Sin especificar'; $data[2,8] = 'Primera'; $data[2,9] = 100; $data[2,10] = 15000; $data[2,11] = 15000; $data[2,12] = 15000; $data[2,13] = '$/caja 16 kilos'; $data[2,14] = 'Provincia de Los Andes'; $data[2,15] = 938; $data[2,16] = 16
$data[3,0] = 44265; $data[3,1] = 9; $data[3,2] = 'Fruta'; $data[3,3] = 100107; $data[3,4] = 'Otros'; $data[3,5] = 100107011; $data[3,6] = 'Tuna'; $data[3,7] = 'Sin especificar'; $data[3,8] = 'Primera'; $data[3,9] = 50; $data[3,10] = 16000; $data[3,11] = 16000; $data[3,12] = 16000; $data[3,13] = '$/caja 16 kilos'; $data[3,14] = 'Provincia de Los Andes'; $data[3,15] = 1000; $data[3,16] = 16
$data[4,0] = 44418; $data[4,1] = 9; $data[4,2] = 'Fruta'; $data[4,3] = 100107; $data[4,4] = 'Otros'; $data[4,5] = 100107011; $data[4,6] = 'Tuna'; $data[4,7] = 'Sin especificar'; $data[4,8] = 'Primera'; $data[4,9] = 40; $data[4,10] = 30000; $data[4,11] = 30000; $data[4,12] = 30000; $data[4,13] = '$/caja 16 kilos'; $data[4,14] = 'Provincia de Los Andes'; $data[4,15] = 1875; $data[4,16] = 16
$data[5,0] = 44270; $data[5,1] = 9; $data[5,2] = 'Fruta'; $data[5,3] = 100107; $data[5,4] = 'Otros'; $data[5,5] = 100107011; $data[5,6] = 'Tuna'; $data[5,7] = 'Sin especificar'; $data[5,8] = 'Primera'; $data[5,9] = 55; $data[5,10] = 10000; $data[5,11] = 10000; $data[5,12] = 10000; $data[5,13] = '$/caja 16 kilos'; $data[5,14] = 'Provincia de Los Andes'; $data[5,15] = 625; $data[5,16] = 16
$data[6,0] = 44630; $data[6,1] = 9; $data[6,2] = 'Fruta'; $data[6,3] = 100107; $data[6,4] = 'Otros'; $data[6,5] = 100107011; $data[6,6] = 'Tuna'; $data[6,7] = 'Sin especificar'; $data[6,8] = 'Primera'; $data[6,9] = 100; $data[6,10] = 17000; $data[6,11] = 17000; $data[6,12] = 17000; $data[6,13] = '$/caja 16 kilos'; $data[6,14] = 'Provincia de Los Andes'; $data[6,15] = 1062; $data[6,16] = 16
$data[7,0] = 44432; $data[7,1] = 9; $data[7,2] = 'Fruta'; $data[7,3] = 100107; $data[7,4] = 'Otros'; $data[7,5] = 100107011; $data[7,6] = 'Tuna'; $data[7,7] = 'Sin especificar'; $data[7,8] = 'Especial'; $data[7,9] = 30; $data[7,10] = 45000; $data[7,11] = 45000; $data[7,12] = 45000; $data[7,13] = '$/caja 18 kilos'; $data[7,14] = 'Provincia de Los Andes'; $data[7,15] = 2500; $data[7,16] = 18
$data[8,0] = 44319; $data[8,1] = 9; $data[8,2] = 'Fruta'; $data[8,3] = 100107; $data[8,4] = 'Otros'; $data[8,5] = 100107011; $data[8,6] = 'Tuna'; $data[8,7] = 'Sin especificar'; $data[8,8] = 'Primera'; $data[8,9] = 50; $data[8,10] = 15000; $data[8,11] = 15000; $data[8,12] = 15000; $data[8,13] = '$/caja 16 kilos'; $data[8,14] = 'Provincia de Los Andes'; $data[8,15] = 938; $data[8,16] = 16
$data[9,0] = 44330; $data[9,1] = 9; $data[9,2] = 'Fruta'; $data[9,3] = 100107; $data[9,4] = 'Otros'; $data[9,5] = 100107011; $data[9,6] = 'Tuna'; $data[9,7] = 'Sin especificar'; $data[9,8] = 'Especial'; $data[9,9] = 35; $data[9,10] = 30000; $data[9,11] = 30000; $data[9,12] = 30000; $data[9,13] = '$/caja 20 kilos'; $data[9,14] = 'Provincia de Los Andes'; $data[9,15] = 1500; $data[9,16] = 20
$data[10,0] = 44326; $data[10,1] = 9; $data[10,2] = 'Fruta'; $data[10,3] = 100107; $data[10,4] = 'Otros'; $data[10,5] = 100107011; $data[10,6] = 'Tuna'; $data[10,7] = 'Sin especificar'; $data[10,8] = 'Primera'; $data[10,9] = 100; $data[10,10] = 15000; $data[10,11] = 15000; $data[10,12] = 15000; $data[10,13] = '$/caja 16 kilos'; $data[10,14] = 'Provincia de Los Andes'; $data[10,15] = 938; $data[10,16] = 16
$data[11,0] = 44253; $data[11,1] = 9; $data[11,2] = 'Fruta'; $data[11,3] = 100107; $data[11,4] = 'Otros'; $data[11,5] = 100107011; $data[11,6] = 'Tuna'; $data[11,7] = 'Sin especificar'; $data[11,8] = 'Primera'; $data[11,9] = 55; $data[11,10] = 16000; $data[11,11] = 16000; $data[11,12] = 16000; $data[11,13] = '$/caja 16 kilos'; $data[11,14] = 'Provincia de Los Andes'; $data[11,15] = 1000; $data[11,16] = 16
$data[12,0] = 44237; $data[12,1] = 9; $data[12,2] = 'Fruta'; $data[12,3] = 100107; $data[12,4] = 'Otros'; $data[12,5] = 100107011; $data[12,6] = 'Tuna'; $data[12,7] = 'Sin especificar'; $data[12,8] = 'Primera'; $data[12,9] = 20; $data[12,10] = 25000; $data[12,11] = 25000; $data[12,12] = 25000; $data[12,13] = '$/caja 16 kilos'; $data[12,14] = 'Provincia de Limarí'; $data[12,15] = 1562; $data[12,16] = 16
$data[13,0] = 44236; $data[13,1] = 9; $data[13,2] = 'Fruta'; $data[13,3] = 100107; $data[13,4] = 'Otros'; $data[13,5] = 100107011; $data[13,6] = 'Tuna'; $data[13,7] = 'Sin especificar'; $data[13,8] = 'Primera'; $data[13,9] = 45; $data[13,10] = 25000; $data[13,11] = 25000; $data[13,12] = 25000; $data[13,13] = '$/caja 16 kilos'; $data[13,14] = 'Provincia de Limarí'; $data[13,15] = 1562; $data[13,16] = 16
$data[14,0] = 44965; $data[14,1] = 9; $data[14,2] = 'Fruta'; $data[14,3] = 100107; $data[14,4] = 'Otros'; $data[14,5] = 100107011; $data[14,6] = 'Tuna'; $data[14,7] = 'Sin especificar'; $data[14,8] = 'Primera'; $data[14,9] = 40; $data[14,10] = 30000; $data[14,11] = 30000; $data[14,12] = 30000; $data[14,13] = '$/caja 18 kilos'; $data[14,14] = 'Provincia de Los Andes'; $data[14,15] = 1667; $data[14,16] = 18
$data[15,0] = 44267; $data[15,1] = 9; $data[15,2] = 'Fruta'; $data[15,3] = 100107; $data[15,4] = 'Otros'; $data[15,5] = 100107011; $data[15,6] = 'Tuna'; $data[15,7] = 'Sin especificar'; $data[15,8] = 'Primera'; $data[15,9] = 50; $data[15,10] = 15000; $data[15,11] = 16000; $data[15,12] = 15400; $data[15,13] = '$/caja 16 kilos'; $data[15,14] = 'Provincia de Los Andes'; $data[15,15] = 962; $data[15,16] = 16
$data[16,0] = 44603; $data[16,1] = 9; $data[16,2] = 'Fruta'; $data[16,3] = 100107; $data[16,4] = 'Otros'; $data[16,5] = 100107011; $data[16,6] = 'Tuna'; $data[16,7] = 'Sin especificar'; $data[16,8] = 'Primera'; $data[16,9] = 40; $data[16,10] = 17000; $data[16,11] = 18000; $data[16,12] = 17500; $data[16,13] = '$/caja 16 kilos'; $data[16,14] = 'Provincia de Los Andes'; $data[16,15] = 1094; $data[16,16] = 16
$data[17,0] = 44323; $data[17,1] = 9; $data[17,2] = 'Fruta'; $data[17,3] = 100107; $data[17,4] = 'Otros'; $data[17,5] = 100107011; $data[17,6] = 'Tuna'; $data[17,7] = 'Sin especificar'; $data[17,8] = 'Segunda'; $data[17,9] = 30; $data[17,10] = 14000; $data[17,11] = 14000; $data[17,12] = 14000; $data[17,13] = '$/caja 16 kilos'; $data[17,14] = 'Provincia de Los Andes'; $data[17,15] = 875; $data[17,16] = 16
$data[18,0] = 44596; $data[18,1] = 9; $data[18,2] = 'Fruta'; $data[18,3] = 100107; $data[18,4] = 'Otros'; $data[18,5] = 100107011; $data[18,6] = 'Tuna'; $data[18,7] = 'Sin especificar'; $data[18,8] = 'Primera'; $data[18,9] = 55; $data[18,10] = 16000; $data[18,11] = 16000; $data[18,12] = 16000; $data[18,13] = '$/caja 16 kilos'; $data[18,14] = 'Provincia de Los Andes'; $data[18,15] = 1000; $data[18,16] = 16
$data[19,0] = 44952; $data[19,1] = 9; $data[19,2] = 'Fruta'; $data[19,3] = 100107; $data[19,4] = 'Otros'; $data[19,5] = 100107011; $data[19,6] = 'Tuna'; $data[19,7] = 'Sin especificar'; $data[19,8] = 'Especial'; $data[19,9] = 50; $data[19,10] = 32000; $data[19,11] = 32000; $data[19,12] = 32000; $data[19,13] = '$/caja 18 kilos'; $data[19,14] = 'Provincia de Los Andes'; $data[19,15] = 1778; $data[19,16] = 18
$data[20,0] = 44435; $data[20,1] = 9; $data[20,2] = 'Fruta'; $data[20,3] = 100107; $data[20,4] = 'Otros'; $data[20,5] = 100107011; $data[20,6] = 'Tuna'; $data[20,7] = 'Sin especificar'; $data[20,8] = 'Especial'; $data[20,9] = 30; $data[20,10] = 45000; $data[20,11] = 45000; $data[20,12] = 45000; $data[20,13] = '$/caja 18 kilos'; $data[20,14] = 'Provincia de Los Andes'; $data[20,15] = 2500; $data[20,16] = 18
$data[21,0] = 44273; $data[21,1] = 9; $data[21,2] = 'Fruta'; $data[21,3] = 100107; $data[21,4] = 'Otros'; $data[21,5] = 100107011; $data[21,6] = 'Tuna'; $data[21,7] = 'Sin especificar'; $data[21,8] = 'Primera'; $data[21,9] = 55; $data[21,10] = 14000; $data[21,11] = 14000; $data[21,12] = 14000; $data[21,13] = '$/caja 16 kilos'; $data[21,14] = 'Provincia de Los Andes'; $data[21,15] = 875; $data[21,16] = 16
$data[22,0] = 44967; $data[22,1] = 9; $data[22,2] = 'Fruta'; $data[22,3] = 100107; $data[22,4] = 'Otros'; $data[22,5] = 100107011; $data[22,6] = 'Tuna'; $data[22,7] = 'Sin especificar'; $data[22,8] = 'Especial'; $data[22,9] = 80; $data[22,10] = 25000; $data[22,11] = 25000; $data[22,12] = 25000; $data[22,13] = '$/caja 18 kilos'; $data[22,14] = 'Provincia de Los Andes'; $data[22,15] = 1389; $data[22,16] = 18
$data[23,0] = 44389; $data[23,1] = 9; $data[23,2] = 'Fruta'; $data[23,3] = 100107; $data[23,4] = 'Otros'; $data[23,5] = 100107011; $data[23,6] = 'Tuna'; $data[23,7] = 'Sin especificar'; $data[23,8] = 'Especial'; $data[23,9] = 30; $data[23,10] = 30000; $data[23,11] = 30000; $data[23,12] = 30000; $data[23,13] = '$/caja 16 kilos'; $data[23,14] = 'Provincia de Los Andes'; $data[23,15] = 1875; $data[23,16] = 16
$data[24,0] = 44594; $data[24,1] = 9; $data[24,2] = 'Fruta'; $data[24,3] = 100107; $data[24,4] = 'Otros'; $data[24,5] = 100107011; $data[24,6] = 'Tuna'; $data[24,7] = 'Sin especificar'; $data[24,8] = 'Primera'; $data[24,9] = 95; $data[24,10] = 16000; $data[24,11] = 16000; $data[24,12] = 16000; $data[24,13] = '$/caja 16 kilos'; $data[24,14] = 'Provincia de Los Andes'; $data[24,15] = 1000; $data[24,16] = 16
$data[25,0] = 44427; $data[25,1] = 9; $data[25,2] = 'Fruta'; $data[25,3] = 100107; $data[25,4] = 'Otros'; $data[25,5] = 100107011; $data[25,6] = 'Tuna'; $data[25,7] = 'Sin especificar'; $data[25,8] = 'Primera'; $data[25,9] = 65; $data[25,10] = 40000; $data[25,11] = 40000; $data[25,12] = 40000; $data[25,13] = '$/caja 16 kilos'; $data[25,14] = 'Provincia de Los Andes'; $data[25,15] = 2500; $data[25,16] = 16
$data[26,0] = 44277; $data[26,1] = 9; $data[26,2] = 'Fruta'; $data[26,3] = 100107; $data[26,4] = 'Otros'; $data[26,5] = 100107011; $data[26,6] = 'Tuna'; $data[26,7] = 'Sin especificar'; $data[26,8] = 'Primera'; $data[26,9] = 90; $data[26,10] = 16000; $data[26,11] = 16000; $data[26,12] = 16000; $data[26,13] = '$/caja 16 kilos'; $data[26,14] = 'Provincia de Los Andes'; $data[26,15] = 1000; $data[26,16] = 16
$data[27,0] = 44221; $data[27,1] = 9; $data[27,2] = 'Fruta'; $data[27,3] = 100107; $data[27,4] = 'Otros'; $data[27,5] = 100107011; $data[27,6] = 'Tuna'; $data[27,7] = 'Sin especificar'; $data[27,8] = 'Primera'; $data[27,9] = 30; $data[27,10] = 25000; $data[27,11] = 25000; $data[27,12] = 25000; $data[27,13] = '$/caja 18 kilos'; $data[27,14] = 'Región Metropolitana'; $data[27,15] = 1389; $data[27,16] = 18
$data[28,0] = 44284; $data[28,1] = 9; $data[28,2] = 'Fruta'; $data[28,3] = 100107; $data[28,4] = 'Otros'; $data[28,5] = 100107011; $data[28,6] = 'Tuna'; $data[28,7] = 'Sin especificar'; $data[28,8] = 'Primera'; $data[28,9] = 25; $data[28,10] = 18000; $data[28,11] = 18000; $data[28,12] = 18000; $data[28,13] = '$/caja 16 kilos'; $data[28,14] = 'Provincia de Los Andes'; $data[28,15] = 1125; $data[28,16] = 16
$data[29,0] = 44708; $data[29,1] = 9; $data[29,2] = 'Fruta'; $data[29,3] = 100107; $data[29,4] = 'Otros'; $data[29,5] = 100107011; $data[29,6] = 'Tuna'; $data[29,7] = 'Sin especificar'; $data[29,8] = 'Primera'; $data[29,9] = 45; $data[29,10] = 19000; $data[29,11] = 19000; $data[29,12] = 19000; $data[29,13] = '$/caja 16 kilos'; $data[29,14] = 'Provincia de Los Andes'; $data[29,15] = 1188; $data[29,16] = 16
$data[30,0] = 44966; $data[30,1] = 9; $data[30,2] = 'Fruta'; $data[30,3] = 100107; $data[30,4] = 'Otros'; $data[30,5] = 100107011; $data[30,6] = 'Tuna'; $data[30,7] = 'Sin especificar'; $data[30,8] = 'Primera'; $data[30,9] = 150; $data[30,10] = 22000; $data[30,11] = 22000; $data[30,12] = 22000; $data[30,13] = '$/caja 18 kilos'; $data[30,14] = 'Provincia de Los Andes'; $data[30,15] = 1222; $data[30,16] = 18
$data[31,0] = 44258; $data[31,1] = 9; $data[31,2] = 'Fruta'; $data[31,3] = 100107; $data[31,4] = 'Otros'; $data[31,5] = 100107011; $data[31,6] = 'Tuna'; $data[31,7] = 'Sin especificar'; $data[31,8] = 'Primera'; $data[31,9] = 65; $data[31,10] = 16000; $data[31,11] = 16000; $data[31,12] = 16000; $data[31,13] = '$/caja 18 kilos granel'; $data[31,14] = 'Provincia de Los Andes'; $data[31,15] = 889; $data[31,16] = 18
$data[32,0] = 44868; $data[32,1] = 9; $data[32,2] = 'Fruta'; $data[32,3] = 100107; $data[32,4] = 'Otros'; $data[32,5] = 100107011; $data[32,6] = 'Tuna'; $data[32,7] = 'Sin especificar'; $data[32,8] = 'Primera'; $data[32,9] = 30; $data[32,10] = 40000; $data[32,11] = 40000; $data[32,12] = 40000; $data[32,13] = '$/caja 16 kilos'; $data[32,14] = 'Provincia de Los Andes'; $data[32,15] = 2500; $data[32,16] = 16
$data[33,0] = 44238; $data[33,1] = 9; $data[33,2] = 'Fruta'; $data[33,3] = 100107; $data[33,4] = 'Otros'; $data[33,5] = 100107011; $data[33,6] = 'Tuna'; $data[33,7] = 'Sin especificar'; $data[33,8] = 'Primera'; $data[33,9] = 65; $data[33,10] = 14000; $data[33,11] = 14000; $data[33,12] = 14000; $data[33,13] = '$/caja 16 kilos'; $data[33,14] = 'Provincia de Los Andes'; $data[33,15] = 875; $data[33,16] = 16
$data[34,0] = 44301; $data[34,1] = 9; $data[34,2] = 'Fruta'; $data[34,3] = 100107; $data[34,4] = 'Otros'; $data[34,5] = 100107011; $data[34,6] = 'Tuna'; $data[34,7] = 'Sin especificar'; $data[34,8] = 'Especial'; $data[34,9] = 55; $data[34,10] = 22000; $data[34,11] = 22000; $data[34,12] = 22000; $data[34,13] = '$/caja 16 kilos'; $data[34,14] = 'Provincia de Los Andes'; $data[34,15] = 1375; $data[34,16] = 16
$data[35,0] = 44301; $data[35,1] = 9; $data[35,2] = 'Fruta'; $data[35,3] = 100107; $data[35,4] = 'Otros'; $data[35,5] = 100107011; $data[35,6] = 'Tuna'; $data[35,7] = 'Sin especificar'; $data[35,8] = 'Primera'; $data[35,9] = 85; $data[35,10] = 19000; $data[35,11] = 19000; $data[35,12] = 19000; $data[35,13] = '$/caja 16 kilos'; $data[35,14] = 'Provincia de Los Andes'; $data[35,15] = 1188; $data[35,16] = 16
$data[36,0] = 44957; $data[36,1] = 9; $data[36,2] = 'Fruta'; $data[36,3] = 100107; $data[36,4] = 'Otros'; $data[36,5] = 100107011; $data[36,6] = 'Tuna'; $data[36,7] = 'Sin especificar'; $data[36,8] = 'Primera'; $data[36,9] = 25; $data[36,10] = 30000; $data[36,11] = 30000; $data[36,12] = 30000; $data[36,13] = '$/caja 18 kilos'; $data[36,14] = 'Provincia de Los Andes'; $data[36,15] = 1667; $data[36,16] = 18
$data[37,0] = 44649; $data[37,1] = 9; $data[37,2] = 'Fruta'; $data[37,3] = 100107; $data[37,4] = 'Otros'; $data[37,5] = 100107011; $data[37,6] = 'Tuna'; $data[37,7] = 'Sin especificar'; $data[37,8] = 'Primera'; $data[37,9] = 65; $data[37,10] = 20000; $data[37,11] = 20000; $data[37,12] = 20000; $data[37,13] = '$/caja 16 kilos'; $data[37,14] = 'Provincia de Los Andes'; $data[37,15] = 1250; $data[37,16] = 16
$data[38,0] = 44649; $data[38,1] = 9; $data[38,2] = 'Fruta'; $data[38,3] = 100107; $data[38,4] = 'Otros'; $data[38,5] = 100107011; $data[38,6] = 'Tuna'; $data[38,7] = 'Sin especificar'; $data[38,8] = 'Segunda'; $data[38,9] = 30; $data[38,10] = 14000; $data[38,11] = 14000; $data[38,12] = 14000; $data[38,13] = '$/caja 16 kilos'; $data[38,14] = 'Provincia de Los Andes'; $data[38,15] = 875; $data[38,16] = 16
$data[39,0] = 44588; $data[39,1] = 9; $data[39,2] = 'Fruta'; $data[39,3] = 100107; $data[39,4] = 'Otros'; $data[39,5] = 100107011; $data[39,6] = 'Tuna'; $data[39,7] = 'Sin especificar'; $data[39,8] = 'Primera'; $data[39,9] = 50; $data[39,10] = 25000; $data[39,11] = 25000; $data[39,12] = 25000; $data[39,13] = '$/caja 16 kilos'; $data[39,14] = 'Provincia de Los Andes'; $data[39,15] = 1562; $data[39,16] = 16
$data[40,0] = 44601; $data[40,1] = 9; $data[40,2] = 'Fruta'; $data[40,3] = 100107; $data[40,4] = 'Otros'; $data[40,5] = 100107011; $data[40,6] = 'Tuna'; $data[40,7] = 'Sin especificar'; $data[40,8] = 'Especial'; $data[40,9] = 30; $data[40,10] = 25000; $data[40,11] = 25000; $data[40,12] = 25000; $data[40,13] = '$/caja 18 kilos'; $data[40,14] = 'Provincia de Los Andes'; $data[40,15] = 1389; $data[40,16] = 18
$data[41,0] = 44601; $data[41,1] = 9; $data[41,2] = 'Fruta'; $data[41,3] = 100107; $data[41,4] = 'Otros'; $data[41,5] = 100107011; $data[41,6] = 'Tuna'; $data[41,7] = 'Sin especificar'; $data[41,8] = 'Primera'; $data[41,9] = 80; $data[41,10] = 18000; $data[41,11] = 18000; $data[41,12] = 18000; $data[41,13] = '$/caja 18 kilos'; $data[41,14] = 'Provincia de Los Andes'; $data[41,15] = 1000; $data[41,16] = 18
$data[42,0] = 44606; $data[42,1] = 9; $data[42,2] = 'Fruta'; $data[42,3] = 100107; $data[42,4] = 'Otros'; $data[42,5] = 100107011; $data[42,6] = 'Tuna'; $data[42,7] = 'Sin especificar'; $data[42,8] = 'Segunda'; $data[42,9] = 80; $data[42,10] = 10000; $data[42,11] = 10000; $data[42,12] = 10000; $data[42,13] = '$/caja 16 kilos'; $data[42,14] = 'Provincia de Los Andes'; $data[42,15] = 625; $data[42,16] = 16
$data[43,0] = 44315; $data[43,1] = 9; $data[43,2] = 'Fruta'; $data[43,3] = 100107; $data[43,4] = 'Otros'; $data[43,5] = 100107011; $data[43,6] = 'Tuna'; $data[43,7] = 'Sin especificar'; $data[43,8] = 'Primera'; $data[43,9] = 105; $data[43,10] = 15000; $data[43,11] = 16000; $data[43,12] = 15619; $data[43,13] = '$/caja 16 kilos'; $data[43,14] = 'Provincia de Los Andes'; $data[43,15] = 976; $data[43,16] = 16
$data[44,0] = 44306; $data[44,1] = 9; $data[44,2] = 'Fruta'; $data[44,3] = 100107; $data[44,4] = 'Otros'; $data[44,5] = 100107011; $data[44,6] = 'Tuna'; $data[44,7] = 'Sin especificar'; $data[44,8] = 'Especial'; $data[44,9] = 50; $data[44,10] = 22000; $data[44,11] = 22000; $data[44,12] = 22000; $data[44,13] = '$/caja 16 kilos'; $data[44,14] = 'Provincia de Los Andes'; $data[44,15] = 1375; $data[44,16] = 16
$data[45,0] = 44294; $data[45,1] = 9; $data[45,2] = 'Fruta'; $data[45,3] = 100107; $data[45,4] = 'Otros'; $data[45,5] = 100107011; $data[45,6] = 'Tuna'; $data[45,7] = 'Sin especificar'; $data[45,8] = 'Primera'; $data[45,9] = 80; $data[45,10] = 14000; $data[45,11] = 16000; $data[45,12] = 15000; $data[45,13] = '$/caja 16 kilos'; $data[45,14] = 'Provincia de Los Andes'; $data[45,15] = 938; $data[45,16] = 16
$data[46,0] = 44244; $data[46,1] = 9; $data[46,2] = 'Fruta'; $data[46,3] = 100107; $data[46,4] = 'Otros'; $data[46,5] = 100107011; $data[46,6] = 'Tuna'; $data[46,7] = 'Sin especificar'; $data[46,8] = 'Primera'; $data[46,9] = 25; $data[46,10] = 14000; $data[46,11] = 14000; $data[46,12] = 14000; $data[46,13] = '$/caja 16 kilos'; $data[46,14] = 'Provincia de Los Andes'; $data[46,15] = 875; $data[46,16] = 16
$data[47,0] = 44280; $data[47,1] = 9; $data[47,2] = 'Fruta'; $data[47,3] = 100107; $data[47,4] = 'Otros'; $data[47,5] = 100107011; $data[47,6] = 'Tuna'; $data[47,7] = 'Sin especificar'; $data[47,8] = 'Especial'; $data[47,9] = 80; $data[47,10] = 20000; $data[47,11] = 20000; $data[47,12] = 20000; $data[47,13] = '$/caja 16 kilos'; $data[47,14] = 'Provincia de Los Andes'; $data[47,15] = 1250; $data[47,16] = 16
$data[48,0] = 44707; $data[48,1] = 9; $data[48,2] = 'Fruta'; $data[48,3] = 100107; $data[48,4] = 'Otros'; $data[48,5] = 100107011; $data[48,6] = 'Tuna'; $data[48,7] = 'Sin especificar'; $data[48,8] = 'Primera'; $data[48,9] = 55; $data[48,10] = 18000; $data[48,11] = 20000; $data[48,12] = 19273; $data[48,13] = '$/caja 16 kilos'; $data[48,14] = 'Provincia de Los Andes'; $data[48,15] = 1205; $data[48,16] = 16
$data[49,0] = 44705; $data[49,1] = 9; $data[49,2] = 'Fruta'; $data[49,3] = 100107; $data[49,4] = 'Otros'; $data[49,5] = 100107011; $data[49,6] = 'Tuna'; $data[49,7] = 'Sin especificar'; $data[49,8] = 'Primera'; $data[49,9] = 25; $data[49,10] = 20000; $data[49,11] = 20000; $data[49,12] = 20000; $data[49,13] = '$/caja 16 kilos'; $data[49,14] = 'Provincia de Los Andes'; $data[49,15] = 1250; $data[49,16] = 16
$data[50,0] = 44964; $data[50,1] = 9; $data[50,2] = 'Fruta'; $data[50,3] = 100107; $data[50,4] = 'Otros'; $data[50,5] = 100107011; $data[50,6] = 'Tuna'; $data[50,7] = 'Sin especificar'; $data[50,8] = 'Especial'; $data[50,9] = 80; $data[50,10] = 30000; $data[50,11] = 30000; $data[50,12] = 30000; $data[50,13] = '$/caja 18 kilos'; $data[50,14] = 'Provincia de Los Andes'; $data[50,15] = 1667; $data[50,16] = 18
$data[51,0] = 44249; $data[51,1] = 9; $data[51,2] = 'Fruta'; $data[51,3] = 100107; $data[51,4] = 'Otros'; $data[51,5] = 100107011; $data[51,6] = 'Tuna'; $data[51,7] = 'Sin especificar'; $data[51,8] = 'Primera'; $data[51,9] = 110; $data[51,10] = 16000; $data[51,11] = 16000; $data[51,12] = 16000; $data[51,13] = '$/caja 16 kilos'; $data[51,14] = 'Provincia de Los Andes'; $data[51,15] = 1000; $data[51,16] = 16
$data[52,0] = 44832; $data[52,1] = 9; $data[52,2] = 'Fruta'; $data[52,3] = 100107; $data[52,4] = 'Otros'; $data[52,5] = 100107011; $data[52,6] = 'Tuna'; $data[52,7] = 'Sin especificar'; $data[52,8] = 'Especial'; $data[52,9] = 50; $data[52,10] = 40000; $data[52,11] = 40000; $data[52,12] = 40000; $data[52,13] = '$/caja 18 kilos'; $data[52,14] = 'Provincia de Los Andes'; $data[52,15] = 2222; $data[52,16] = 18
$data[53,0] = 44274; $data[53,1] = 9; $data[53,2] = 'Fruta'; $data[53,3] = 100107; $data[53,4] = 'Otros'; $data[53,5] = 100107011; $data[53,6] = 'Tuna'; $data[53,7] = 'Sin especificar'; $data[53,8] = 'Especial'; $data[53,9] = 35; $data[53,10] = 16000; $data[53,11] = 16000; $data[53,12] = 16000; $data[53,13] = '$/caja 16 kilos'; $data[53,14] = 'Provincia de Los Andes'; $data[53,15] = 1000; $data[53,16] = 16
$data[54,0] = 44274; $data[54,1] = 9; $data[54,2] = 'Fruta'; $data[54,3] = 100107; $data[54,4] = 'Otros'; $data[54,5] = 100107011; $data[54,6] = 'Tuna'; $data[54,7] = 'Sin especificar'; $data[54,8] = 'Primera'; $data[54,9] = 40; $data[54,10] = 14000; $data[54,11] = 14000; $data[54,12] = 14000; $data[54,13] = '$/caja 16 kilos'; $data[54,14] = 'Provincia de Los Andes'; $data[54,15] = 875; $data[54,16] = 16
$data[55,0] = 44663; $data[55,1] = 9; $data[55,2] = 'Fruta'; $data[55,3] = 100107; $data[55,4] = 'Otros'; $data[55,5] = 100107011; $data[55,6] = 'Tuna'; $data[55,7] = 'Sin especificar'; $data[55,8] = 'Primera'; $data[55,9] = 55; $data[55,10] = 16000; $data[55,11] = 16000; $data[55,12] = 16000; $data[55,13] = '$/caja 16 kilos'; $data[55,14] = 'Provincia de Los Andes'; $data[55,15] = 1000; $data[55,16] = 16
$data[56,0] = 44266; $data[56,1] = 9; $data[56,2] = 'Fruta'; $data[56,3] = 100107; $data[56,4] = 'Otros'; $data[56,5] = 100107011; $data[56,6] = 'Tuna'; $data[56,7] = 'Sin especificar'; $data[56,8] = 'Especial'; $data[56,9] = 50; $data[56,10] = 18000; $data[56,11] = 19000; $data[56,12] = 18400; $data[56,13] = '$/caja 16 kilos'; $data[56,14] = 'Provincia de Los Andes'; $data[56,15] = 1150; $data[56,16] = 16
$data[57,0] = 44266; $data[57,1] = 9; $data[57,2] = 'Fruta'; $data[57,3] = 100107; $data[57,4] = 'Otros'; $data[57,5] = 100107011; $data[57,6] = 'Tuna'; $data[57,7] = 'Sin especificar'; $data[57,8] = 'Primera'; $data[57,9] = 120; $data[57,10] = 15000; $data[57,11] = 15000; $data[57,12] = 15000; $data[57,13] = '$/caja 16 kilos'; $data[57,14] = 'Provincia de Los Andes'; $data[57,15] = 938; $data[57,16] = 16
$data[58,0] = 44266; $data[58,1] = 9; $data[58,2] = 'Fruta'; $data[58,3] = 100107; $data[58,4] = 'Otros'; $data[58,5] = 100107011; $data[58,6] = 'Tuna'; $data[58,7] = 'Sin especificar'; $data[58,8] = 'Segunda'; $data[58,9] = 20; $data[58,10] = 10000; $data[58,11] = 10000; $data[58,12] = 10000; $data[58,13] = '$/caja 16 kilos'; $data[58,14] = 'Provincia de Los Andes'; $data[58,15] = 625; $data[58,16] = 16
$data[59,0] = 44658; $data[59,1] = 9; $data[59,2] = 'Fruta'; $data[59,3] = 100107; $data[59,4] = 'Otros'; $data[59,5] = 100107011; $data[59,6] = 'Tuna'; $data[59,7] = 'Sin especificar'; $data[59,8] = 'Primera'; $data[59,9] = 200; $data[59,10] = 17000; $data[59,11] = 17000; $data[59,12] = 17000; $data[59,13] = '$/caja 16 kilos'; $data[59,14] = 'Provincia de Los Andes'; $data[59,15] = 1062; $data[59,16] = 16
$data[60,0] = 44637; $data[60,1] = 9; $data[60,2] = 'Fruta'; $data[60,3] = 100107; $data[60,4] = 'Otros'; $data[60,5] = 100107011; $data[60,6] = 'Tuna'; $data[60,7] = 'Sin especificar'; $data[60,8] = 'Primera'; $data[60,9] = 65; $data[60,10] = 15000; $data[60,11] = 15000; $data[60,12] = 15000; $data[60,13] = '$/caja 16 kilos'; $data[60,14] = 'Provincia de Los Andes'; $data[60,15] = 938; $data[60,16] = 16
$data[61,0] = 44650; $data[61,1] = 9; $data[61,2] = 'Fruta'; $data[61,3] = 100107; $data[61,4] = 'Otros'; $data[61,5] = 100107011; $data[61,6] = 'Tuna'; $data[61,7] = 'Sin especificar'; $data[61,8] = 'Primera'; $data[61,9] = 45; $data[61,10] = 18000; $data[61,11] = 20000; $data[61,12] = 19333; $data[61,13] = '$/caja 16 kilos'; $data[61,14] = 'Provincia de Los Andes'; $data[61,15] = 1208; $data[61,16] = 16
$data[62,0] = 44235; $data[62,1] = 9; $data[62,2] = 'Fruta'; $data[62,3] = 100107; $data[62,4] = 'Otros'; $data[62,5] = 100107011; $data[62,6] = 'Tuna'; $data[62,7] = 'Sin especificar'; $data[62,8] = 'Primera'; $data[62,9] = 55; $data[62,10] = 25000; $data[62,11] = 25000; $data[62,12] = 25000; $data[62,13] = '$/caja 16 kilos'; $data[62,14] = 'Provincia de Limarí'; $data[62,15] = 1562; $data[62,16] = 16
$data[63,0] = 44320; $data[63,1] = 9; $data[63,2] = 'Fruta'; $data[63,3] = 100107; $data[63,4] = 'Otros'; $data[63,5] = 100107011; $data[63,6] = 'Tuna'; $data[63,7] = 'Sin especificar'; $data[63,8] = 'Primera'; $data[63,9] = 50; $data[63,10] = 14000; $data[63,11] = 14000; $data[63,12] = 14000; $data[63,13] = '$/caja 16 kilos'; $data[63,14] = 'Provincia de Los Andes'; $data[63,15] = 875; $data[63,16] = 16
$data[64,0] = 44251; $data[64,1] = 9; $data[64,2] = 'Fruta'; $data[64,3] = 100107; $data[64,4] = 'Otros'; $data[64,5] = 100107011; $data[64,6] = 'Tuna'; $data[64,7] = 'Sin especificar'; $data[64,8] = 'Primera'; $data[64,9] = 55; $data[64,10] = 16000; $data[64,11] = 16000; $data[64,12] = 16000; $data[64,13] = '$/caja 16 kilos'; $data[64,14] = 'Provincia de Los Andes'; $data[64,15] = 1000; $data[64,16] = 16
$data[65,0] = 44595; $data[65,1] = 9; $data[65,2] = 'Fruta'; $data[65,3] = 100107; $data[65,4] = 'Otros'; $data[65,5] = 100107011; $data[65,6] = 'Tuna'; $data[65,7] = 'Sin especificar'; $data[65,8] = 'Primera'; $data[65,9] = 115; $data[65,10] = 16000; $data[65,11] = 16000; $data[65,12] = 16000; $data[65,13] = '$/caja 16 kilos'; $data[65,14] = 'Provincia de Los Andes'; $data[65,15] = 1000; $data[65,16] = 16

# Columns A:C are identical ('10' / 'Vega Modelo de Temuco' / 'La Araucanía')
# for every data row, including the two new rows created at the bottom.
$abc = New-Object 'object[,]' 2,3
$abc[0,0] = 10; $abc[0,1] = 'Vega Modelo de Temuco'; $abc[0,2] = 'La Araucanía'
$abc[1,0] = 10; $abc[1,1] = 'Vega Modelo de Temuco'; $abc[1,2] = 'La Araucanía'
$ws.Range("A71:C72").Value2 = $abc

# Write the D:T block and re-apply the date number format used by column D
# (Value2 writes raw date serials; the format must be (re)applied to the two
# newly created rows so they render as dates just like the existing ones).
$ws.Range("D7:T72").Value2 = $data
$ws.Range("D7:D72").NumberFormat = "YYYY-MM-DD HH:MM:SS"
